# Updated video recording code to add proper synchronization
$wb = $excel.ActiveWorkbook

# --- "Count" sheet: mark the invalid-credentials test case to Execute (TRUE) ---
$wsCount = $wb.Worksheets.Item("Count")
$wsCount.Range("B3").Value = $true

# --- "TestData" sheet: enable execution and switch run mode to remote ---
$wsData = $wb.Worksheets.Item("TestData")
$wsData.Range("F5").Value = $true
$wsData.Range("G5").Value = "remote"

# Update the active selection to reflect the last edited cell
$wsData.Activate()
$wsData.Range("G5").Select()
